# Fix bug impacting Problem 2 results: several search methods actually
# timed out on problem 2, so their row data must be replaced with a
# "timeout" marker, and the rows that DID complete need their corrected
# numbers. Mirrors the "Problem 3" table's existing timeout-row layout
# (rows 30/32/34/38, which use the "timeout" shared string with the
# italic/grey style already present in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TimeoutRow([int]$row) {
    # Clear the "heuristic" column (D), mark Plan Length (E) as the
    # existing "timeout" text style (copy format from a known-good
    # timeout row), and blank out the remaining numeric columns.
    $ws.Range("D$row").Value = $null

    $ws.Range("E30").Copy() | Out-Null
    $ws.Range("E$row").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("E$row").Value = "timeout"

    $ws.Range("F$row").Value = $null
    $ws.Range("G$row").Value = $null
    $ws.Range("H$row").Value = $null
    $ws.Range("I$row").Value = $null
}

# Row 16: breadth_first_search - corrected completed numbers
$ws.Range("D16").Value = $null
$ws.Range("E16").Value = 9
$ws.Range("F16").Value = 3343
$ws.Range("G16").Value = 4609
$ws.Range("H16").Value = 30509
$ws.Range("I16").Value = 16.798993903997999

# Row 17: breadth_first_tree_search - now a timeout
Set-TimeoutRow 17

# Row 18: depth_first_graph_search - corrected completed numbers
$ws.Range("E18").Value = 575
$ws.Range("F18").Value = 582
$ws.Range("G18").Value = 583
$ws.Range("H18").Value = 5211
$ws.Range("I18").Value = 3.72759502197732

# Row 19: depth_limited_search - now a timeout
Set-TimeoutRow 19

# Row 20: uniform_cost_search - corrected completed numbers
$ws.Range("D20").Value = $null
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 4761
$ws.Range("G20").Value = 4763
$ws.Range("H20").Value = 43206
$ws.Range("I20").Value = 15.701451928995001

# Row 21: recursive_best_first_search h_1 - now a timeout
Set-TimeoutRow 21

# Row 22: greedy_best_first_graph_search h_1 - corrected completed numbers
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 550
$ws.Range("G22").Value = 552
$ws.Range("H22").Value = 4950
$ws.Range("I22").Value = 1.73293168400414

# Row 23: astar_search h_1 - corrected completed numbers
$ws.Range("D23").Value = $null
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 4761
$ws.Range("G23").Value = 4763
$ws.Range("H23").Value = 43206
$ws.Range("I23").Value = 15.489678598009

# Row 24: astar_search h_ignore_preconditions - corrected completed
# numbers, and un-bold the row (it's no longer the last real result in
# the table now that rows below it are timeouts, so it loses the
# "final row" emphasis style).
$ws.Range("C24:I24").Font.Bold = $false
$ws.Range("D24").Borders.LineStyle = -4142   # xlLineStyleNone
$ws.Range("D24").IndentLevel = 0
$ws.Range("D24").HorizontalAlignment = -4131 # xlLeft
$ws.Range("D24").Value = $null
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 1450
$ws.Range("G24").Value = 1452
$ws.Range("H24").Value = 13303
$ws.Range("I24").Value = 4.6805308510083696

# Row 25: astar_search h_pg_levelsum - now a timeout
Set-TimeoutRow 25

Write-Output "Problem 2 results corrected"
